# Add a new (blank-looking) paragraph containing three spaces right after
# the "- Luu: VT, HSCB." paragraph in the "Noi nhan:" table cell, matching
# the formatting (color/size/lang) of the paragraph it follows.

$d = $word.ActiveDocument

# Locate the anchor paragraph robustly via Find rather than a hard-coded
# paragraph index.
$rng = $d.Content
$found = $rng.Find.Execute("- Lưu: VT, HSCB.", $false, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0)

if ($found) {
    # Insert right at the end of the matched text (i.e. just before the
    # paragraph mark that ends "- Luu: VT, HSCB."). Inserting a full <w:p>
    # there splits it off as a new paragraph that keeps the original
    # paragraph's mark and cell, which is what keeps the new paragraph
    # inside the same table cell instead of spilling into the next cell.
    $insertPoint = $d.Range($rng.End, $rng.End)

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:pPr>' + `
             '<w:rPr>' + `
               '<w:color w:val="000000" w:themeColor="text1"/>' + `
               '<w:sz w:val="22"/>' + `
               '<w:szCs w:val="24"/>' + `
               '<w:lang w:val="nl-NL"/>' + `
             '</w:rPr>' + `
           '</w:pPr>' + `
           '<w:r>' + `
             '<w:rPr>' + `
               '<w:color w:val="000000" w:themeColor="text1"/>' + `
               '<w:sz w:val="22"/>' + `
               '<w:szCs w:val="24"/>' + `
               '<w:lang w:val="nl-NL"/>' + `
             '</w:rPr>' + `
             '<w:t xml:space="preserve">   </w:t>' + `
           '</w:r>' + `
           '</w:p>'

    [void]$insertPoint.InsertXML($xml)
}
